# Sprint 1 presentation edit:
# Slide 5, Title placeholder ("Title 1") text update.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Hier zijn we nog niet aan toegekomen"
